$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values - repulled data / mean calculation fix
$updates = @{
    3  = -6
    4  = -5
    7  = -2
    11 = -4
    13 = -2
    15 = -12
    16 = 4
    24 = -4
    26 = -1
    27 = 1
    28 = 1
    30 = 1
    31 = -2
    36 = 3
    41 = 0
    47 = -3
    48 = 5
    52 = 1
    53 = -1
    54 = 0
    55 = -4
    57 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
